$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.080.87'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.90%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.888.15'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.68%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9990'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.15%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.21'
$ws.Range('D5').Style = 'Normal'

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9993'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.10%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5140'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.07%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3725'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.61%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07210'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.71%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9037'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.39%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.05'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.88%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07639'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.37%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.903.35'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.63%  '

$ws.Range('E14').Value = '  +2.05%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.263'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.81%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9993'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.11%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008506'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.04%  '

$ws.Range('E18').Value = '  +2.32%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9988'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.17%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '27.092.94'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.78%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.059'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.96%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.130.86'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.10%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.56'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.49%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.418'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.09%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.30'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.90%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.794'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.03%  '

$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.180'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +6.33%  '

$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.04'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.25%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.975'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.91%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.837'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.32%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09216'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.08%  '

$ws.Range('E33').Value = '  -0.35%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7686'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.26%  '

$ws.Range('E35').Value = '  +4.18%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.971'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.01%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.277'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.87%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.585'
$ws.Range('D38').Style = 'Normal'

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5641'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.03%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01995'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.58%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.079'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.23%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.939'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.04%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '118.79'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.09%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.591'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.66%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1500'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.69%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4817'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.89%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9991'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.08%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.13'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.22%  '

$ws.Range('E49').Value = '  +2.08%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '37.21'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.99%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '64.19'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.40%  '
